# Update the "Short Name" for row 4 (K4) to include the bowl size.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K4").Value = 'Gerber Maxwell ADA EL Bowl 14"'

# Widen column K to fit the new, longer short-name text (drop the old
# "best fit" auto-sizing in favor of an explicit custom width).
$ws.Columns.Item(11).ColumnWidth = 45.665

# Scroll the view over and select K2:K5 (the Short Name column) as the
# new active selection.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("K2:K5").Select() | Out-Null

# Set page setup to letter-ish paper size 9 (A4) / portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
